$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.525.87"
$ws.Range("E2").Value = "  +1.77%  "
$ws.Range("D3").Value = "1.567.82"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  -1.54%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.00"
$ws.Range("E5").Value = "  +1.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.489"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.990"
$ws.Range("E7").Value = "  -1.58%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.66"
$ws.Range("E8").Value = "  +2.64%  "
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0595"
$ws.Range("E10").Value = "  -0.34%  "
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").Value = "1.565.42"
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("D16").Value = "27.516.51"
$ws.Range("E16").Value = "  +1.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.41"
$ws.Range("E17").Value = "  +0.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "225.44"
$ws.Range("E18").Value = "  +4.16%  "
$ws.Range("E19").Value = "  +1.54%  "
$ws.Range("D20").Value = "0.0₃0705"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("E21").Value = "  -1.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.12"
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.44"
$ws.Range("E23").Value = "  +2.42%  "
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.83"
$ws.Range("E25").Value = "  -2.66%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.62"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.17"
$ws.Range("E27").Value = "  +0.62%  "
$ws.Range("E28").Value = "  +1.51%  "
$ws.Range("E29").Value = "  -1.47%  "
$ws.Range("E30").Value = "  +1.21%  "
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("D33").Value = "1.448.56"
$ws.Range("E33").Value = "  +1.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.15"
$ws.Range("E34").Value = "  -1.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.12"
$ws.Range("E35").Value = "  +2.20%  "
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("E39").Value = "  +1.56%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.72"
$ws.Range("E41").Value = "  -1.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.35"
$ws.Range("E42").Value = "  +0.94%  "
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("E44").Value = "  +5.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.973"
$ws.Range("E45").Value = "  -3.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.57"
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("D47").Value = "1.704.23"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.68"
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0526"
$ws.Range("E49").Value = "  +1.48%  "
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.990"
$ws.Range("E51").Value = "  -1.51%  "
